# Reorder item rows within each brand group (Item Name / UOM / TP travel together)
# plus update the cumulative sales target (TP) figures, per commit:
# "update dashboard.png, banner_ai.png and cumilative sales target"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Dinafex 120mg <-> Dinafex 180mg Tablet (row 3 / row 4)
$ws.Range("C3").Value = "Dinafex 180mg Tablet"
$ws.Range("BB3").Value = 224.89

$ws.Range("C4").Value = "Dinafex 120mg Tablet"
$ws.Range("BB4").Value = 179.91

# Etorix 120mg / 60mg-40's / 90mg rotate (row 7 / row 8 / row 9)
$ws.Range("C7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D7").Value = "40's"

$ws.Range("C8").Value = "Etorix 90mg Tablet"
$ws.Range("D8").Value = "30's"
$ws.Range("BB8").Value = 269.87

$ws.Range("C9").Value = "Etorix 120mg Tablet"
$ws.Range("D9").Value = "20's"
$ws.Range("BB9").Value = 209.9

# Flucloxin 500mg Capsule <-> Flucloxin 500mg Capsule - 36's (row 11 / row 12)
$ws.Range("C11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D11").Value = "36 's"
$ws.Range("BB11").Value = 284.21

$ws.Range("C12").Value = "Flucloxin 500mg Capsule"
$ws.Range("D12").Value = "30 's"
$ws.Range("BB12").Value = 237.74

# Ketonic 10mg Tablet <-> Ketonic 30mg Injection (row 15 / row 16)
$ws.Range("C15").Value = "Ketonic 30mg Injection"
$ws.Range("D15").Value = "5 's"
$ws.Range("BB15").Value = 206.77

$ws.Range("C16").Value = "Ketonic 10mg Tablet"
$ws.Range("D16").Value = "20's"
$ws.Range("BB16").Value = 150.38

# Kynol TR 200mg Capsule <-> Kynol TR 100mg Capsule (row 18 / row 19)
$ws.Range("C18").Value = "Kynol TR 100mg Capsule"
$ws.Range("D18").Value = "50 's"
$ws.Range("BB18").Value = 262.37

$ws.Range("C19").Value = "Kynol TR 200mg Capsule"
$ws.Range("D19").Value = "30 's"
$ws.Range("BB19").Value = 224.89

# Zithrox 500mg Tablet / 30ml Dry Suspension / 250mg Tablet-6's rotate (row 26 / row 27 / row 28)
$ws.Range("C26").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("D26").Value = "6's"
$ws.Range("BB26").Value = 89.95999999999999

$ws.Range("C27").Value = "Zithrox 500mg Tablet"
$ws.Range("D27").Value = "6 's"
$ws.Range("BB27").Value = 136.83

$ws.Range("C28").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D28").Value = "30ml"
$ws.Range("BB28").Value = 97.45
